# A new weekly price record for "Perejil" (Vega Central Mapocho de Santiago)
# is inserted as the new row 187. Every existing row from 187 downward
# (old 187..308) shifts down by one (to 188..309); the workbook's used
# range grows from A1:R308 to A1:R309.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 187..308 down to 188..309, carrying their values/formatting.
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with this week's record.
$ws.Cells.Item(187, 1).Value  = 9
$ws.Cells.Item(187, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(187, 3).Value  = "Metropolitana"
$ws.Cells.Item(187, 4).Value  = 44603
$ws.Cells.Item(187, 5).Value  = 13
$ws.Cells.Item(187, 6).Value  = 100112044
$ws.Cells.Item(187, 7).Value  = "Perejil"
$ws.Cells.Item(187, 8).Value  = "Sin especificar"
$ws.Cells.Item(187, 9).Value  = "Primera"
$ws.Cells.Item(187, 10).Value = 61
$ws.Cells.Item(187, 11).Value = 16000
$ws.Cells.Item(187, 12).Value = 18000
$ws.Cells.Item(187, 13).Value = 17016
$ws.Cells.Item(187, 14).Value = "$/docena de atados"
$ws.Cells.Item(187, 15).Value = "Región Metropolitana"
$ws.Cells.Item(187, 16).Value = 5672
$ws.Cells.Item(187, 17).Value = 3
$ws.Cells.Item(187, 18).Value = "Hortaliza"
